# Update Betfair Back/Lay odds cells to reflect the latest refresh of the
# Jogos_do_Dia sheet (commit: "Atualizando o arquivo XLSX").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (ACS Petrolul 52 vs Unirea Slobozia) ---
$ws.Range("I2").Value = 6.2
$ws.Range("J2").Value = 3.15
$ws.Range("L2").Value = 1.01
$ws.Range("M2").Value = 1.09
$ws.Range("N2").Value = 2.64
$ws.Range("O2").Value = 1.49
$ws.Range("R2").Value = 1.2
$ws.Range("S2").Value = 5
$ws.Range("T2").Value = 1.98
$ws.Range("U2").Value = 1.6
$ws.Range("V2").Value = 1.2
$ws.Range("W2").Value = 2
$ws.Range("X2").Value = 9.6
$ws.Range("Y2").Value = 14.5
$ws.Range("Z2").Value = 44
$ws.Range("AA2").Value = 210
$ws.Range("AB2").Value = 6.8
$ws.Range("AC2").Value = 8.199999999999999
$ws.Range("AD2").Value = 24
$ws.Range("AE2").Value = 130
$ws.Range("AF2").Value = 10.5
$ws.Range("AG2").Value = 11.5
$ws.Range("AH2").Value = 27
$ws.Range("AI2").Value = 140
$ws.Range("AJ2").Value = 24
$ws.Range("AK2").Value = 27
$ws.Range("AL2").Value = 60
$ws.Range("AM2").Value = 260
$ws.Range("AN2").Value = 22
$ws.Range("AO2").Value = 210

# --- Row 4 (Universitatea Cluj vs Arges Pitesti) ---
$ws.Range("F4").Value = 2.34
$ws.Range("G4").Value = 2.82
$ws.Range("K4").Value = 3.25

# --- Row 7 (Al-Ettifaq vs Al-Taawoun Buraidah) ---
$ws.Range("F7").Value = 2.46
$ws.Range("G7").Value = 2.62
$ws.Range("I7").Value = 3.05
$ws.Range("J7").Value = 3.7
$ws.Range("K7").Value = 4.2
$ws.Range("P7").Value = 2.1
$ws.Range("Q7").Value = 1.76

# --- Row 8 (Hermannstadt vs Rapid Bucharest) ---
$ws.Range("J8").Value = 3.4
$ws.Range("P8").Value = 1.77
$ws.Range("Q8").Value = 2.08

# --- Row 9 (FC Dordrecht vs Helmond Sport) ---
$ws.Range("F9").Value = 1.82
$ws.Range("G9").Value = 1.99
$ws.Range("H9").Value = 4
$ws.Range("I9").Value = 4.7
$ws.Range("J9").Value = 3.9
$ws.Range("K9").Value = 4.7
$ws.Range("P9").Value = 2.38
$ws.Range("Q9").Value = 1.61

# --- Row 13 (Altrincham vs Morecambe) ---
$ws.Range("F13").Value = 1.91
$ws.Range("G13").Value = 2.28
$ws.Range("K13").Value = 5

# --- Row 14 (Rochdale vs Boreham Wood) ---
$ws.Range("F14").Value = 1.94
$ws.Range("G14").Value = 2.04

# --- Row 15 (Bologna vs AC Milan) ---
$ws.Range("F15").Value = 3.65
$ws.Range("G15").Value = 3.75
$ws.Range("H15").Value = 2.24
$ws.Range("I15").Value = 2.28
$ws.Range("J15").Value = 3.45
$ws.Range("L15").Value = 1.42
$ws.Range("P15").Value = 1.86
$ws.Range("V15").Value = 1.78
$ws.Range("W15").Value = 1.36
$ws.Range("Y15").Value = 9.4
$ws.Range("Z15").Value = 14
$ws.Range("AA15").Value = 30
$ws.Range("AB15").Value = 13
$ws.Range("AC15").Value = 7.8
$ws.Range("AD15").Value = 11.5
$ws.Range("AE15").Value = 26
$ws.Range("AF15").Value = 25
$ws.Range("AG15").Value = 15.5
$ws.Range("AJ15").Value = 70
$ws.Range("AK15").Value = 44
$ws.Range("AL15").Value = 55
$ws.Range("AN15").Value = 46
$ws.Range("AO15").Value = 19.5

# --- Row 16 (St Mirren vs Hearts) ---
$ws.Range("P16").Value = 1.82
